# BIS-1002: Fixed XLS export tests
# Adds a new "Internal Assignment" column (O) to the experiment type export sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a piece of text into a cell without Excel auto-converting
# strings such as "TRUE"/"FALSE" into boolean cells. We do this by writing
# a formula that evaluates to the literal text, then collapsing the formula
# down to its value via copy / paste-special (values only). This preserves
# the cell's existing style (number format / font) untouched.
function Set-Text($addr, [string]$text) {
    $r = $ws.Range($addr)
    $escaped = $text -replace '"', '""'
    $r.Formula = '="' + $escaped + '"'
    $r.Copy() | Out-Null
    $r.PasteSpecial(-4163) | Out-Null
}

# Start from a clean slate (values only - formatting/styles stay intact) so
# that the shared string table gets rebuilt in the exact order the cells
# are (re)populated below, matching the order cells are laid out on sheet.
$ws.Cells.ClearContents() | Out-Null

# Give the new column O a header style matching the existing bold header
# cells (K4:N4) but bumped up to size 12 - this creates the new font +
# cell style used for the "Internal Assignment" header.
$ws.Range("K4").Copy() | Out-Null
$ws.Range("O4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("O4").Font.Size = 12

# Row 1
Set-Text "A1" "EXPERIMENT_TYPE"

# Row 2
Set-Text "A2" "Description"
Set-Text "B2" "Code"

# Row 3
Set-Text "A3" "Default experiment"
Set-Text "B3" "DEFAULT_EXPERIMENT"

# Row 4 (header row) - includes the new "Internal Assignment" column O
Set-Text "A4" "Code"
Set-Text "B4" "Mandatory"
Set-Text "C4" "Show in edit views"
Set-Text "D4" "Section"
Set-Text "E4" "Property label"
Set-Text "F4" "Data type"
Set-Text "G4" "Vocabulary code"
Set-Text "H4" "Description"
Set-Text "I4" "Metadata"
Set-Text "J4" "Dynamic script"
Set-Text "K4" "Multivalued"
Set-Text "L4" "Unique"
Set-Text "M4" "Pattern"
Set-Text "N4" "Pattern Type"
Set-Text "O4" "Internal Assignment"

# Row 5
Set-Text "A5" "`$NAME"
Set-Text "B5" "FALSE"
Set-Text "C5" "TRUE"
Set-Text "D5" "General info"
Set-Text "E5" "Name"
Set-Text "F5" "VARCHAR"
Set-Text "H5" "Name"
Set-Text "K5" "FALSE"
Set-Text "L5" "FALSE"
Set-Text "O5" "FALSE"

# Row 6
Set-Text "A6" "`$DEFAULT_OBJECT_TYPE"
Set-Text "B6" "FALSE"
Set-Text "C6" "TRUE"
Set-Text "D6" "General info"
Set-Text "E6" "Default object type"
Set-Text "F6" "VARCHAR"
Set-Text "H6" "Enter the code of the object type for which the collection is used"
Set-Text "K6" "FALSE"
Set-Text "L6" "FALSE"
Set-Text "O6" "FALSE"

# Row 7
Set-Text "A7" "NOTES"
Set-Text "B7" "FALSE"
Set-Text "C7" "TRUE"
Set-Text "E7" "Notes"
Set-Text "F7" "MULTILINE_VARCHAR"
Set-Text "H7" "Notes"
Set-Text "I7" '{"custom_widget":"Word Processor"}'
Set-Text "K7" "FALSE"
Set-Text "L7" "FALSE"
Set-Text "O7" "FALSE"

# Row 8
Set-Text "A8" "`$XMLCOMMENTS"
Set-Text "B8" "FALSE"
Set-Text "C8" "FALSE"
Set-Text "E8" "Comments List"
Set-Text "F8" "XML"
Set-Text "H8" "Comments log"
Set-Text "K8" "FALSE"
Set-Text "L8" "FALSE"
Set-Text "O8" "FALSE"

$excel.CutCopyMode = $false

# Match the author's final selection (O7:O8) as seen in the saved workbook.
$ws.Range("O7:O8").Select() | Out-Null
